$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fill in F/G/H for existing rows 182-186 (previously blank) ---
$ws.Cells.Item(182, 6).Value = 0
$ws.Cells.Item(182, 7).Value = 24
$ws.Cells.Item(182, 8).Value = 4
$ws.Cells.Item(183, 6).Value = 1
$ws.Cells.Item(183, 7).Value = 24
$ws.Cells.Item(183, 8).Value = 4
$ws.Cells.Item(184, 6).Value = 2
$ws.Cells.Item(184, 7).Value = 24
$ws.Cells.Item(184, 8).Value = 4
$ws.Cells.Item(185, 6).Value = 3
$ws.Cells.Item(185, 7).Value = 64
$ws.Cells.Item(185, 8).Value = 5
$ws.Cells.Item(186, 6).Value = 4
$ws.Cells.Item(186, 7).Value = 64
$ws.Cells.Item(186, 8).Value = 5

# --- Append new experiment rows 187-199 ---
# Row 187
$ws.Cells.Item(187, 1).Value = "Beto-19032020-001"
$ws.Cells.Item(187, 2).Value = "200319_Beto_rfMapper_basic"
$ws.Cells.Item(187, 3).Value = "N:\Stimuli\2019-06-RF-mapping\2020-03-19-Beto"
$row187_D = @"
001 RFMapping 13:
-8:8:8 
Carlos' huge image test
Completed
"@
$ws.Cells.Item(187, 4).Value = $row187_D
$ws.Cells.Item(187, 5).Value = "ReducDimen_Evol"

# Row 188
$ws.Cells.Item(188, 1).Value = "Beto-19032020-002"
$ws.Cells.Item(188, 2).Value = "200319_Beto_rfMapper_basic(1)"
$ws.Cells.Item(188, 3).Value = "N:\Stimuli\2019-06-RF-mapping\2020-03-19-Beto"
$row188_D = @"
002 RFMapping starts 13:36
-8:2:8
Mid image, 6mins
Completed
"@
$ws.Cells.Item(188, 4).Value = $row188_D
$ws.Cells.Item(188, 5).Value = "ReducDimen_Evol"

# Row 189
$ws.Cells.Item(189, 1).Value = "Beto-19032020-003"
$ws.Cells.Item(189, 2).Value = "200319_Beto_rfMapper_basic(2)"
$ws.Cells.Item(189, 3).Value = "N:\Stimuli\2019-06-RF-mapping\2020-03-19-Beto"
$row189_D = @"
003 RFMapping starts 13:49
-2:1:2
Completed
"@
$ws.Cells.Item(189, 4).Value = $row189_D
$ws.Cells.Item(189, 5).Value = "ReducDimen_Evol"

# Row 190
$ws.Cells.Item(190, 1).Value = "Beto-19032020-004"
$ws.Cells.Item(190, 2).Value = "200319_Beto_generate_integrated"
$ws.Cells.Item(190, 3).Value = "N:\Stimuli\2019-12-Evolutions\2020-03-19-Beto-01\2020-03-19-14-06-29"
$row190_D = @"
004 Generate Integrated starts 14:06
25 [-1 -1.5 ] 4 1 ZOHA_Sphere_lr_euclid
25 [-1 -1.5 ] 4 1 ZOHA_Sphere_lr_euclid_RD
A little bit single unit, Maybe we should not split this unit.
Still use the Inverse decay + 1.2 learning rate protocol
Response is sparse.
Starts taking off around gen 7, very slowly.
Not very success to take off
Not very successful in both……Not an informative comparison?
Gen 13, Exploration 41.8 deg seems to take off…… Not really just flucturation.
It's a good sign, Full evolution starts to dominate around 20 gens. Around 25 gens reduced starts to take off as well.
Interestingly this evolution is pretty unstable, the score keep going up and down.
Seems they plateau at around very different spots.
Still growing! Exploration 29.2 - 28.2 degs. (This cell is more sensititve to finer changes? )
Pretty large gap!
27.3 deg, the cell still growng… Maybe single units really likes small exploration? The interaction goes.
Really taking off! Around 45!
Seems you really need to have a small learning rate for this cell to get to the top of mountain right?
How could we test the sharpness of mountain hypothesis? Successfulness of different learning rate tuning schedule? 
Great evolution! 51gens, the ZOHA full finally plateaued. 51 gens
Completed
"@
$ws.Cells.Item(190, 4).Value = $row190_D
$ws.Cells.Item(190, 5).Value = "ReducDimen_Evol"

# Row 191
$ws.Cells.Item(191, 1).Value = "Beto-19032020-005"
$ws.Cells.Item(191, 2).Value = "200319_Beto_generate_integrated(1)"
$ws.Cells.Item(191, 3).Value = "N:\Stimuli\2019-12-Evolutions\2020-03-19-Beto-02\2020-03-19-15-07-34"
$row191_D = @"
005 Generate Integrated starts 15:07
25 [-1 -1.5 ] 4 1 CMAES
See if CMAES can replicate the really late growing of response.
Just curious, not relavent to current exp series, but can serve as control for the final activation level and the learning curve.
CMAES goes up and down pretty fast.
Seems CMAES's learning curve is smoother.
Starts bumpping around.
Seems it plateaus pretty early and didn't get the peak up there! Ends 60gens
Completed
"@
$ws.Cells.Item(191, 4).Value = $row191_D
$ws.Cells.Item(191, 5).Value = "ReducDimen_Evol"

# Row 192
$ws.Cells.Item(192, 1).Value = "Beto-19032020-006"
$ws.Cells.Item(192, 2).Value = "200319_Beto_generate_integrated(2)"
$ws.Cells.Item(192, 3).Value = "N:\Stimuli\2019-12-Evolutions\2020-03-19-Beto-03\2020-03-19-15-50-58"
$row192_D = @"
006 Generate Integrated starts 15:50
39 [0 0 ] 3 2 ZOHA_Sphere_lr_euclid
39 [0 0 ] 3 2 ZOHA_Sphere_lr_euclid_RD
A fast V1 Hash Evolution comparison.
Comes pretty close, just as predicted.
maybe 11 generations are enough. Bump him up. 16mins 15 gens
Completed
"@
$ws.Cells.Item(192, 4).Value = $row192_D
$ws.Cells.Item(192, 5).Value = "ReducDimen_Evol"

# Row 193
$ws.Cells.Item(193, 1).Value = "Beto-20032020-001"
$ws.Cells.Item(193, 2).Value = "200320_Beto_rfMapper_basic"
$ws.Cells.Item(193, 3).Value = "N:\Stimuli\2019-06-RF-mapping\2020-03-20-Beto"
$row193_D = @"
001 rf mapping -8:8:8
Completed
"@
$ws.Cells.Item(193, 4).Value = $row193_D
$ws.Cells.Item(193, 5).Value = "ReducDimen_Evol"

# Row 194
$ws.Cells.Item(194, 1).Value = "Beto-20032020-002"
$ws.Cells.Item(194, 2).Value = "200320_Beto_rfMapper_basic(1)"
$ws.Cells.Item(194, 3).Value = "N:\Stimuli\2019-06-RF-mapping\2020-03-20-Beto"
$row194_D = @"
002 rf mapping -4:4:4j
Completed
"@
$ws.Cells.Item(194, 4).Value = $row194_D
$ws.Cells.Item(194, 5).Value = "ReducDimen_Evol"

# Row 195
$ws.Cells.Item(195, 1).Value = "Beto-20032020-003"
$ws.Cells.Item(195, 2).Value = "200320_Beto_rfMapper_basic(2)"
$ws.Cells.Item(195, 3).Value = "N:\Stimuli\2019-06-RF-mapping\2020-03-20-Beto"
$row195_D = @"
003 rf mapping -2:2:2
Completed
"@
$ws.Cells.Item(195, 4).Value = $row195_D
$ws.Cells.Item(195, 5).Value = "ReducDimen_Evol"

# Row 196
$ws.Cells.Item(196, 1).Value = "Beto-20032020-004"
$ws.Cells.Item(196, 2).Value = "200320_Beto_rfMapper_basic(3)"
$ws.Cells.Item(196, 3).Value = "N:\Stimuli\2019-06-RF-mapping\2020-03-20-Beto"
$row196_D = @"
004 rf mapping -4:2:4
Completed
"@
$ws.Cells.Item(196, 4).Value = $row196_D
$ws.Cells.Item(196, 5).Value = "ReducDimen_Evol"

# Row 197
$ws.Cells.Item(197, 1).Value = "Beto-20032020-005"
$ws.Cells.Item(197, 2).Value = "200320_Beto_generate_integrated(1)"
$ws.Cells.Item(197, 3).Value = "N:\Stimuli\2019-12-Evolutions\2020-03-20-Beto-01\2020-03-20-09-41-52"
$row197_D = @"
005 at 941 AM
integrated
5 (0,0) 4 1 ZOHA full
5 (0,0) 4 1 ZOHA red
finished, 40 gens
The BHV files seem to be deleted at first so I recover it and it's correct now!
Basis saved.
Completed
"@
$ws.Cells.Item(197, 4).Value = $row197_D
$ws.Cells.Item(197, 5).Value = "ReducDimen_Evol"

# Row 198
$ws.Cells.Item(198, 1).Value = "Beto-20032020-006"
$ws.Cells.Item(198, 2).Value = "200320_Beto_generate_integrated"
$ws.Cells.Item(198, 3).Value = "N:\Stimuli\2019-12-Evolutions\2020-03-20-Beto-02\2020-03-20-10-36-13"
$row198_D = @"
006 at 1036
33 (0,0) 3 1 zoha full
33 (0,0) 3 1 zoha red
1046 AM complaining
31 blocks completed…
complaining
BHV name is strange but it's coV name is strange burrect!
Basis saved.
Completed
"@
$ws.Cells.Item(198, 4).Value = $row198_D
$ws.Cells.Item(198, 5).Value = "ReducDimen_Evol"

# Row 199
$ws.Cells.Item(199, 1).Value = "Beto-20032020-007"
$ws.Cells.Item(199, 2).Value = "200320_Beto_generate_integrated(3)"
$ws.Cells.Item(199, 3).Value = "N:\Stimuli\2019-12-Evolutions\2020-03-20-Beto-03\2020-03-20-11-17-46"
$row199_D = @"
007 1117
58 (0,0) 3 1  zoha full
58 (0,0) 3 1 zoha red
Blasting him, last evolution for the day
Bumping 1148
Basis saved.
Completed
"@
$ws.Cells.Item(199, 4).Value = $row199_D

